$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from E1 into F1
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

# Populate time_taken column with per-row timestamps as text
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:39:39.590631"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:39:39.590642"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:39:39.590645"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:39:39.590647"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:39:39.590650"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:39:39.590653"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:39:39.590655"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:39:39.590658"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:39:39.590660"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:39:39.590663"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:39:39.590665"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:39:39.590668"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:39:39.590670"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:39:39.590673"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:39:39.590675"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:39:39.590678"
